$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.5085426078324898, 0.6312770643337786, 0.5775225091701024, 0.7599490174808455, 0.5845391272395405, 15)
    3  = @(0.3604583747834282, 0.5284892497187151, 0.3749043735392947, 0.6122943520393559, 0.5136324539103049, 14)
    4  = @(0.2671644164584084, 0.4711099050880848, 0.3212841183425387, 0.5668192995501641, 0.5203199983321719, 13)
    5  = @(0.4282641160228051, 0.548481490052795, 0.4014706193385897, 0.6336170920505457, 0.4877336824715215, 12)
    6  = @(0.4297788858055521, 0.5712726212628527, 0.4264280183742889, 0.6530145621456607, 0.5156451691415619, 11)
    7  = @(0.3783814472866451, 0.5371069045676157, 0.3902711452809265, 0.6247168520865485, 0.5239790343878171, 10)
    8  = @(0.3326355726653664, 0.4945997455591428, 0.3233547881867544, 0.5686429355815075, 0.4891798334788348, 9)
    9  = @(0.3620304170102688, 0.5407351962130211, 0.3743586110673328, 0.6118485197067431, 0.5273032071390139, 8)
    10 = @(0.4970330926187166, 0.5512696337016082, 0.4018017306810739, 0.6338783248235216, 0.4249154521118755, 7)
    11 = @(0.4078060000532727, 0.4615113944401893, 0.259491323550082, 0.5094029088551439, 0.3343990248860763, 6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
